$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1123.0132
$ws.Range("I15").Value = 1123.0132
$ws.Range("K15").Value = 3369.0396
$ws.Range("M15").Value = -3200.0396

$ws.Range("H100").Value = 2333.3333
$ws.Range("I100").Value = 1800
$ws.Range("K100").Value = 1800
$ws.Range("M100").Value = -1259

$ws.Range("H137").Value = 5715393
$ws.Range("I137").Value = 1015.86365
$ws.Range("J137").Value = 15385877
$ws.Range("K137").Value = 3047.59095
$ws.Range("L137").Value = 46157631
$ws.Range("M137").Value = -497.5909499999998
$ws.Range("N137").Value = -46162731

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14603.31
$ws.Range("I32").Value = 16147.547
$ws.Range("J32").Value = 1734.6666
$ws.Range("K32").Value = 16147.547
$ws.Range("L32").Value = 1734.6666
$ws.Range("M32").Value = -15860.547
$ws.Range("N32").Value = -2308.6666

$ws.Range("H45").Value = 1000
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()

$ws.Range("H74").Value = 1800
$ws.Range("I74").Value = 1800
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1800
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -926
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 1800
$ws.Range("I77").Value = 1800
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 9000
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -4632
$ws.Range("N77").ClearContents()

$ws.Range("H122").Value = 6983.8945
$ws.Range("I122").Value = 1667.75
$ws.Range("K122").Value = 5003.25
$ws.Range("M122").Value = -2553.25

$ws.Range("H132").Value = 3803.7068
$ws.Range("I132").Value = 4285.909
$ws.Range("J132").Value = 3167.2
$ws.Range("K132").Value = 12857.727
$ws.Range("L132").Value = 9501.599999999999
$ws.Range("M132").Value = -10327.727
$ws.Range("N132").Value = -14561.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1626.6216
$ws.Range("I105").Value = 1360.4166
$ws.Range("J105").Value = 2118.077
$ws.Range("K105").Value = 1360.4166
$ws.Range("L105").Value = 2118.077
$ws.Range("M105").Value = 386.5834
$ws.Range("N105").Value = -5612.077

$ws.Range("H134").Value = 19573.482
$ws.Range("I134").Value = 25940.244
$ws.Range("J134").Value = 2171
$ws.Range("K134").Value = 77820.73199999999
$ws.Range("L134").Value = 6513
$ws.Range("M134").Value = -75285.73199999999
$ws.Range("N134").Value = -11583

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2986683.2
$ws.Range("I31").Value = 1808.7097
$ws.Range("J31").Value = 5556992
$ws.Range("K31").Value = 1808.7097
$ws.Range("L31").Value = 5556992
$ws.Range("M31").Value = -1513.7097
$ws.Range("N31").Value = -5557582

$ws.Range("H34").Value = 2986683.2
$ws.Range("I34").Value = 1808.7097
$ws.Range("J34").Value = 5556992
$ws.Range("K34").Value = 1808.7097
$ws.Range("L34").Value = 5556992
$ws.Range("M34").Value = -1606.7097
$ws.Range("N34").Value = -5557396

$ws.Range("H86").Value = 47621204
$ws.Range("I86").Value = 66667616
$ws.Range("J86").Value = 5168.3335
$ws.Range("K86").Value = 66667616
$ws.Range("L86").Value = 5168.3335
$ws.Range("M86").Value = -66666493
$ws.Range("N86").Value = -7414.3335

$ws.Range("H89").Value = 47621204
$ws.Range("I89").Value = 66667616
$ws.Range("J89").Value = 5168.3335
$ws.Range("K89").Value = 333338080
$ws.Range("L89").Value = 25841.6675
$ws.Range("M89").Value = -333332464
$ws.Range("N89").Value = -37073.6675

$ws.Range("H99").Value = 2079.92
$ws.Range("I99").Value = 1593.625
$ws.Range("J99").Value = 2944.4443
$ws.Range("K99").Value = 1593.625
$ws.Range("L99").Value = 2944.4443
$ws.Range("M99").Value = -95.625
$ws.Range("N99").Value = -5940.4443

$ws.Range("H105").Value = 758.63635
$ws.Range("I105").Value = 543.61536
$ws.Range("J105").Value = 1069.2222
$ws.Range("K105").Value = 543.61536
$ws.Range("L105").Value = 1069.2222
$ws.Range("M105").Value = 1203.38464
$ws.Range("N105").Value = -4563.2222

$ws.Range("H126").Value = 2079.92
$ws.Range("I126").Value = 1593.625
$ws.Range("J126").Value = 2944.4443
$ws.Range("K126").Value = 4780.875
$ws.Range("L126").Value = 8833.332900000001
$ws.Range("M126").Value = -2310.875
$ws.Range("N126").Value = -13773.3329

$ws.Range("H134").Value = 1243.9166
$ws.Range("I134").Value = 1243.9166
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3731.7498
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -1196.7498
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 849.1667
$ws.Range("I15").Value = 100
$ws.Range("J15").Value = 999
$ws.Range("K15").Value = 300
$ws.Range("L15").Value = 2997
$ws.Range("M15").Value = -160
$ws.Range("N15").Value = -3277

$ws.Range("H68").Value = 1327.9241
$ws.Range("I68").Value = 1208.8864
$ws.Range("J68").Value = 1477.5714
$ws.Range("K68").Value = 3626.6592
$ws.Range("L68").Value = 4432.7142
$ws.Range("M68").Value = -2815.6592
$ws.Range("N68").Value = -6054.7142

$ws.Range("H71").Value = 1327.9241
$ws.Range("I71").Value = 1208.8864
$ws.Range("J71").Value = 1477.5714
$ws.Range("K71").Value = 10879.9776
$ws.Range("L71").Value = 13298.1426
$ws.Range("M71").Value = -6823.9776
$ws.Range("N71").Value = -21410.1426

$ws.Range("H107").Value = 881.8946999999999
$ws.Range("I107").Value = 660
$ws.Range("J107").Value = 961.1429000000001
$ws.Range("K107").Value = 1980
$ws.Range("L107").Value = 2883.4287
$ws.Range("M107").Value = -60
$ws.Range("N107").Value = -6723.4287

$ws.Range("H113").Value = 598.8
$ws.Range("J113").Value = 669.625
$ws.Range("L113").Value = 2008.875
$ws.Range("N113").Value = -6348.875

$ws.Range("H122").Value = 947.0476
$ws.Range("I122").Value = 808.46155
$ws.Range("J122").Value = 1172.25
$ws.Range("K122").Value = 7276.15395
$ws.Range("L122").Value = 10550.25
$ws.Range("M122").Value = -4826.15395
$ws.Range("N122").Value = -15450.25

$ws.Range("H131").Value = 5377482.5
$ws.Range("J131").Value = 8217702.5
$ws.Range("L131").Value = 24653107.5
$ws.Range("N131").Value = -24663187.5

$ws.Range("H132").Value = 1774.2
$ws.Range("I132").Value = 876
$ws.Range("J132").Value = 2100.818
$ws.Range("K132").Value = 7884
$ws.Range("L132").Value = 18907.362
$ws.Range("M132").Value = -5354
$ws.Range("N132").Value = -23967.362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I107").Value = 530.1429000000001
$ws.Range("J107").Value = 575.75
$ws.Range("K107").Value = 530.1429000000001
$ws.Range("L107").Value = 575.75
$ws.Range("M107").Value = 1389.8571
$ws.Range("N107").Value = -4415.75

$ws.Range("H122").Value = 8922.223
$ws.Range("I122").Value = 56300
$ws.Range("K122").Value = 168900
$ws.Range("M122").Value = -166450
